$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four occurrence rows (18-21) for this locality got re-matched to a
# different set of GPS points / record Ids. Net effect: row 18 <-> row 21
# swap their Id/Antal-Stadium/Metod/coordinate data, and row 19 <-> row 20
# swap their Id/coordinate data.

# --- Row 18 <-> Row 21 ---
$ws.Range("A18").Value = 111821926
$ws.Range("K18").Value = ""
$ws.Range("Q18").Value = 550846.2444635418
$ws.Range("R18").Value = 6681625.195240833

$ws.Range("A21").Value = 111821923
$ws.Range("K21").Value = "blomning"
$ws.Range("Q21").Value = 550701.1291094749
$ws.Range("R21").Value = 6681909.496304798

# --- Row 19 <-> Row 20 ---
$ws.Range("A19").Value = 111821928
$ws.Range("Q19").Value = 550825.9503372401
$ws.Range("R19").Value = 6681726.144349095

$ws.Range("A20").Value = 111821927
$ws.Range("Q20").Value = 550819.8901872271
$ws.Range("R20").Value = 6681733.007140613

Write-Output "edits applied"
